{"js": "// The placeholder \"{{vendor_name}}\" in the \"BETWEEN: / {{vendor_name}} /\n// (the \"Vendor\")\" paragraph was originally typed as five separate runs\n// (\"{{\", \"vendor\", \"_\", \"name\", \"}}\"). Collapse it back into a single run\n// containing the literal text \"{{vendor_name}}\" (same visible text, same\n// formatting), matching how the other {{vendor_name}} placeholders later\n// in the document already appear as one run.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the specific paragraph (\"BETWEEN:\\u000b{{vendor_name}}\\u000b(the\n// \"Vendor\")\") that still has the placeholder split across runs, instead of\n// assuming a fixed index.\nlet targetParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text.indexOf(\"BETWEEN:\") !== -1 && text.indexOf(\"{{vendor_name}}\") !== -1) {\n    targetParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (targetParagraph) {\n  // Search within that paragraph only, so the other already-intact\n  // \"{{vendor_name}}\" occurrences elsewhere in the document are untouched.\n  const matches = targetParagraph.search(\"{{vendor_name}}\", { matchCase: true });\n  matches.load(\"items\");\n  await context.sync();\n\n  if (matches.items.length > 0) {\n    // Replacing the matched range merges the multiple runs it spans into a\n    // single run carrying the full literal text, keeping the formatting of\n    // the first run in the match (same rFonts/sz/szCs as before).\n    matches.items[0].insertText(\"{{vendor_name}}\", \"Replace\");\n    await context.sync();\n  }\n}\n", "ps1": "# The placeholder \"{{vendor_name}}\" in the \"BETWEEN: / {{vendor_name}} /\n# (the \"Vendor\")\" paragraph was originally typed as five separate runs\n# (\"{{\", \"vendor\", \"_\", \"name\", \"}}\"). Collapse it back into a single run\n# containing the literal text \"{{vendor_name}}\" (same visible text, same\n# formatting), matching how the other {{vendor_name}} placeholders later\n# in the document already appear as one run.\n\n$d = $word.ActiveDocument\n\n# Locate the specific paragraph (\"BETWEEN:\" followed by \"{{vendor_name}}\")\n# that still has the placeholder split across runs, instead of assuming a\n# fixed paragraph index.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -like \"*BETWEEN:*\" -and $t -like \"*{{vendor_name}}*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    # Restrict Find/Replace to this paragraph's range only, so the other\n    # already-intact \"{{vendor_name}}\" occurrences elsewhere in the\n    # document are left completely untouched.\n    $rng = $target.Range\n    $find = $rng.Find\n    $find.Text = \"{{vendor_name}}\"\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Forward = $true\n    $find.Wrap = 0\n\n    # Replacing (even with identical visible text) merges the runs the\n    # match spans into a single run, carrying the formatting of the first\n    # run in the match (same rFonts/sz/szCs as before). wdReplaceOne (1)\n    # replaces only this single found occurrence.\n    $find.Execute(\n        [ref]$find.Text,\n        [ref]$find.MatchCase,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $false,\n        \"{{vendor_name}}\",\n        1\n    ) | Out-Null\n}\n"}
